$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.005445415939153612
$ws.Range("E2").Value = 0.0017732871330551343
$ws.Range("F2").Value = 0.17073597764897733
$ws.Range("G2").Value = 0.14762383638314816
$ws.Range("H2").Value = 0.34769365815003167
$ws.Range("I2").Value = 62.56
$ws.Range("J2").Value = 0.07239999999999995

$ws.Range("D3").Value = 0.00034202999640824735
$ws.Range("E3").Value = 0.0031990425105815645
$ws.Range("F3").Value = 0.05594593704940671
$ws.Range("G3").Value = 0.03904919325530441
$ws.Range("H3").Value = 0.19144788117686348
$ws.Range("I3").Value = 69.98
$ws.Range("J3").Value = 0.11739999999999995

$ws.Range("D4").Value = 0.00035189360940551096
$ws.Range("E4").Value = -0.00018375397850538738
$ws.Range("F4").Value = 0.002122044909120078
$ws.Range("G4").Value = 0.00002847715544966025
$ws.Range("H4").Value = 0.004946043652206112
$ws.Range("I4").Value = 80.0
$ws.Range("J4").Value = 0.31619999999999865

$ws.Range("D5").Value = -0.00005361577160693303
$ws.Range("E5").Value = -0.00001577306097741054
$ws.Range("F5").Value = 0.000819661675654828
$ws.Range("G5").Value = 0.0000018318548429481967
$ws.Range("H5").Value = 0.0010879720792912203
$ws.Range("I5").Value = 77.98
$ws.Range("J5").Value = 0.27800000000000097

$ws.Range("D6").Value = 0.00003424445810232461
$ws.Range("E6").Value = -0.00001016659700736072
$ws.Range("F6").Value = 0.0004102396503581929
$ws.Range("G6").Value = 0.0000012454322841375446
$ws.Range("H6").Value = 0.001048388322574096
$ws.Range("I6").Value = 80.72
$ws.Range("J6").Value = 0.4517999999999995

